$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(27).Insert()

$ws.Cells.Item(27, 1).Value = 7
$ws.Cells.Item(27, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(27, 3).Value = "Ñuble"
$ws.Cells.Item(27, 4).Value = 44547
$ws.Cells.Item(27, 5).Value = 16
$ws.Cells.Item(27, 6).Value = 100112017
$ws.Cells.Item(27, 7).Value = "Apio"
$ws.Cells.Item(27, 8).Value = "Americana (o)"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 60
$ws.Cells.Item(27, 11).Value = 8000
$ws.Cells.Item(27, 12).Value = 8500
$ws.Cells.Item(27, 13).Value = 8250
$ws.Cells.Item(27, 14).Value = "`$/docena de matas"
$ws.Cells.Item(27, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(27, 16).Value = 1375
$ws.Cells.Item(27, 17).Value = 6
$ws.Cells.Item(27, 18).Value = "Hortaliza"
